{"js": "// Remove the trailing \"Ver no Jupiter...\" line, the copyright/footer line,\n// and the blank paragraph that separated them from the preceding\n// \"Requisitos\" block (LOQ4044 requirement line stays).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two text paragraphs that must be removed by their content.\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIndex === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (copyrightIndex === -1 && t.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the paragraphs targeted for removal.\");\n}\n\n// The blank paragraph immediately preceding the \"Ver no Jupiter...\" line\n// (right after the LOQ4044 requirement paragraph) is removed too.\nconst blankIndex = jupiterIndex - 1;\nif (blankIndex < 0 || items[blankIndex].text !== \"\") {\n  throw new Error(\"Unexpected document shape: blank separator paragraph not found.\");\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[jupiterIndex].delete();\nitems[blankIndex].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" line, the copyright/footer line,\n# and the blank paragraph that separated them from the preceding\n# \"Requisitos\" block (the LOQ4044 requirement line itself stays).\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$n = $paras.Count\n\n$jupiterIdx = -1\n$copyrightIdx = -1\nfor ($i = 1; $i -le $n; $i++) {\n    $t = $paras.Item($i).Range.Text\n    if ($jupiterIdx -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIdx = $i\n    }\n    if ($copyrightIdx -eq -1 -and $t -like \"*Powered by Jekyll*\") {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -eq -1 -or $copyrightIdx -eq -1) {\n    throw \"Could not locate the paragraphs targeted for removal.\"\n}\n\n# The blank paragraph immediately preceding the \"Ver no Jupiter...\" line.\n$blankIdx = $jupiterIdx - 1\n$blankText = $paras.Item($blankIdx).Range.Text.Trim()\nif ($blankIdx -lt 1 -or $blankText -ne \"\") {\n    throw \"Unexpected document shape: blank separator paragraph not found.\"\n}\n\n# Delete bottom-up so earlier paragraph indices remain valid.\n$paras.Item($copyrightIdx).Range.Delete()\n$paras.Item($jupiterIdx).Range.Delete()\n$paras.Item($blankIdx).Range.Delete()\n"}
